$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "출퇴근2, 멸망한 파이썬, 3개월 재?직"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222796834109"

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D45").Value = "Anomaly detection - 1-SVM, SVDD"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/108"

$ws.Range("D50").Value = "Kullback-Leibler divergence"
$ws.Range("E50").Value = "http://incredible.egloos.com/7543794"

$ws.Range("D51").Value = "[python] 숫자 1부터 100까지의 리스트 생성하기, range 객체를 list 객체로 변환하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EC%88%AB%EC%9E%90-1%EB%B6%80%ED%84%B0-100%EA%B9%8C%EC%A7%80%EC%9D%98-%EB%A6%AC%EC%8A%A4%ED%8A%B8-%EC%83%9D%EC%84%B1%ED%95%98%EA%B8%B0"
